# Auto-generated edit script: updates the cryptos worksheet
# with refreshed prices / 1h volume %, matching the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value2 = '59.723.11'
$ws.Range('E2').Value2 = '  -6.86%  '

# Row 3
$ws.Range('D3').Value2 = '3.295.95'
$ws.Range('E3').Value2 = '  -5.24%  '

# Row 4
$ws.Range('D4').Value2 = "'0.999"
$ws.Range('E4').Value2 = '  -0.02%  '

# Row 5
$ws.Range('D5').Value2 = "'559.83"
$ws.Range('E5').Value2 = '  -4.33%  '

# Row 6
$ws.Range('D6').Value2 = "'127.56"
$ws.Range('E6').Value2 = '  -2.99%  '

# Row 7
$ws.Range('E7').Value2 = '  -0.07%  '

# Row 8
$ws.Range('D8').Value2 = '3.295.10'
$ws.Range('E8').Value2 = '  -5.23%  '

# Row 9
$ws.Range('E9').Value2 = '  -2.77%  '

# Row 10
$ws.Range('D10').Value2 = "'7.40"
$ws.Range('E10').Value2 = '  -4.09%  '

# Row 11
$ws.Range('E11').Value2 = '  -6.02%  '

# Row 12
$ws.Range('D12').Value2 = "'0.369"
$ws.Range('E12').Value2 = '  -4.41%  '

# Row 13
$ws.Range('D13').Value2 = '3.856.86'
$ws.Range('E13').Value2 = '  -5.22%  '

# Row 14
$ws.Range('E14').Value2 = '  -0.26%  '

# Row 15
$ws.Range('D15').Value2 = '3.298.14'
$ws.Range('E15').Value2 = '  -5.16%  '

# Row 16
$ws.Range('E16').Value2 = '  -7.00%  '

# Row 17
$ws.Range('D17').Value2 = '59.936.42'
$ws.Range('E17').Value2 = '  -6.49%  '

# Row 18
$ws.Range('D18').Value2 = "'23.84"
$ws.Range('E18').Value2 = '  -5.21%  '

# Row 19
$ws.Range('E19').Value2 = '  -1.75%  '

# Row 20
$ws.Range('D20').Value2 = "'13.17"
$ws.Range('E20').Value2 = '  -1.72%  '

# Row 21
$ws.Range('D21').Value2 = "'8.85"
$ws.Range('E21').Value2 = '  -11.33%  '

# Row 22
$ws.Range('D22').Value2 = "'349.72"
$ws.Range('E22').Value2 = '  -9.27%  '

# Row 23
$ws.Range('E23').Value2 = '  -3.11%  '

# Row 24
$ws.Range('E24').Value2 = '  -0.10%  '

# Row 25
$ws.Range('D25').Value2 = '3.423.09'
$ws.Range('E25').Value2 = '  -5.36%  '

# Row 26
$ws.Range('D26').Value2 = "'68.25"
$ws.Range('E26').Value2 = '  -8.36%  '

# Row 27
$ws.Range('E27').Value2 = '  -4.21%  '

# Row 28
$ws.Range('D28').Value2 = "'0.999"
$ws.Range('E28').Value2 = '  -0.23%  '

# Row 29
$ws.Range('E29').Value2 = '  +2.46%  '

# Row 30
$ws.Range('E30').Value2 = '  +0.04%  '

# Row 31
$ws.Range('E31').Value2 = '  -2.66%  '

# Row 32
$ws.Range('E32').Value2 = '  -3.15%  '

# Row 33
$ws.Range('E33').Value2 = '  -6.38%  '

# Row 34
$ws.Range('E34').Value2 = '  +0.00%  '

# Row 35
$ws.Range('D35').Value2 = '3.319.31'
$ws.Range('E35').Value2 = '  -5.32%  '

# Row 36
$ws.Range('D36').Value2 = "'22.61"
$ws.Range('E36').Value2 = '  -1.64%  '

# Row 37
$ws.Range('D37').Value2 = "'5.24"
$ws.Range('E37').Value2 = '  +0.65%  '

# Row 38
$ws.Range('E38').Value2 = '  -0.83%  '

# Row 39
$ws.Range('E39').Value2 = '  -2.20%  '

# Row 40
$ws.Range('D40').Value2 = "'156.46"
$ws.Range('E40').Value2 = '  -4.29%  '

# Row 41
$ws.Range('D41').Value2 = "'0.0743"
$ws.Range('E41').Value2 = '  -4.40%  '

# Row 42
$ws.Range('E42').Value2 = '  +0.19%  '

# Row 43
$ws.Range('D43').Value2 = "'40.32"
$ws.Range('E43').Value2 = '  -3.00%  '

# Row 44
$ws.Range('B44').Value2 = 'Mantle'
$ws.Range('C44').Value2 = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D44').Value2 = "'0.738"
$ws.Range('E44').Value2 = '  -7.52%  '

# Row 45
$ws.Range('B45').Value2 = 'Filecoin'
$ws.Range('C45').Value2 = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D45').Value2 = "'4.27"
$ws.Range('E45').Value2 = '  -1.65%  '

# Row 46
$ws.Range('E46').Value2 = '  +2.45%  '

# Row 47
$ws.Range('B47').Value2 = 'Stacks'
$ws.Range('C47').Value2 = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D47').Value2 = "'1.53"
$ws.Range('E47').Value2 = '  -5.57%  '

# Row 48
$ws.Range('B48').Value2 = 'EnergySwap'
$ws.Range('C48').Value2 = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D48').Value2 = "'22.32"
$ws.Range('E48').Value2 = '  -4.82%  '

# Row 49
$ws.Range('E49').Value2 = '  -0.74%  '

# Row 50
$ws.Range('D50').Value2 = "'21.70"
$ws.Range('E50').Value2 = '  +5.95%  '

# Row 51
$ws.Range('D51').Value2 = "'0.844"
$ws.Range('E51').Value2 = '  -6.07%  '

